# Update COVID-19 country statistics (refresh timestamp 14:40 -> 15:57)
# and re-sort side effects on the "Pais" sheet (sorted descending by "Casos totales").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh timestamp in the title cell (A1)
$ws.Range("A1").Value = 'Datos actualizados a 13 de Julio de 2020 a las 15:57'

# --- Straightforward per-country statistic updates ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 3415649
$ws.Range("C4").Value = 1654
$ws.Range("D4").Value = 1517570
$ws.Range("E4").Value = 1760282
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 137797

# India (row 6)
$ws.Range("B6").Value = 888944
$ws.Range("C6").Value = 9478
$ws.Range("D6").Value = 560307
$ws.Range("E6").Value = 305304
$ws.Range("G6").Value = 146
$ws.Range("H6").Value = 23333

# Arabia Saudita (row 17)
$ws.Range("B17").Value = 235111
$ws.Range("C17").Value = 2852
$ws.Range("D17").Value = 169842
$ws.Range("E17").Value = 63026
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 2243

# Alemania (row 19)
$ws.Range("B19").Value = 200047
$ws.Range("C19").Value = 97
$ws.Range("E19").Value = 5812

# Argentina (row 25)
$ws.Range("D25").Value = 44173
$ws.Range("E25").Value = 54134
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 1859

# Irak (row 28)
$ws.Range("B28").Value = 79735
$ws.Range("C28").Value = 2229
$ws.Range("D28").Value = 46998
$ws.Range("E28").Value = 29487
$ws.Range("G28").Value = 100
$ws.Range("H28").Value = 3250

# Filipinas (row 36)
$ws.Range("C36").Value = 747
$ws.Range("G36").Value = 65

# Paises Bajos (row 40)
$ws.Range("B40").Value = 51093
$ws.Range("C40").Value = 71

# Portugal (row 42)
$ws.Range("B42").Value = 46818
$ws.Range("C42").Value = 306
$ws.Range("D42").Value = 31065
$ws.Range("E42").Value = 14091
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 1662

# Israel (row 46)
$ws.Range("B46").Value = 39871
$ws.Range("C46").Value = 1201
$ws.Range("D46").Value = 19256
$ws.Range("E46").Value = 20251
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 364

# Serbia (row 63)
$ws.Range("B63").Value = 18639
$ws.Range("C63").Value = 279
$ws.Range("D63").Value = 13940
$ws.Range("E63").Value = 4294
$ws.Range("G63").Value = 12
$ws.Range("H63").Value = 405

# Republica de Macedonia (row 81)
$ws.Range("B81").Value = 8197
$ws.Range("C81").Value = 86
$ws.Range("D81").Value = 4326
$ws.Range("E81").Value = 3486
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 385

# Bosnia y Herzegovina (row 87)
$ws.Range("B87").Value = 6981
$ws.Range("C87").Value = 104
$ws.Range("D87").Value = 3179
$ws.Range("E87").Value = 3576
$ws.Range("G87").Value = 5
$ws.Range("H87").Value = 226

# Mali (row 112)
$ws.Range("B112").Value = 2412
$ws.Range("C112").Value = 1
$ws.Range("D112").Value = 1730
$ws.Range("E112").Value = 561

# Tunez (row 133)
$ws.Range("B133").Value = 1302
$ws.Range("C133").Value = 39
$ws.Range("D133").Value = 1082
$ws.Range("E133").Value = 170

# --- Rows whose case counts caused them to swap sort position ---
# Azerbaiyan <-> Ghana (rows 57/58)
$ws.Range("A57").Value = 'Azerbaiyan'
$ws.Range("B57").Value = 24570
$ws.Range("C57").Value = 529
$ws.Range("D57").Value = 15640
$ws.Range("E57").Value = 8617
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 313
$ws.Range("A58").Value = 'Ghana'
$ws.Range("B58").Value = 24518
$ws.Range("D58").Value = 20187
$ws.Range("E58").Value = 4192
$ws.Range("H58").Value = 139

# Etiopia <-> Costa Rica (rows 83/84)
$ws.Range("A83").Value = 'Etiopia'
$ws.Range("B83").Value = 7766
$ws.Range("C83").Value = 206
$ws.Range("D83").Value = 2430
$ws.Range("E83").Value = 5208
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 128
$ws.Range("A84").Value = 'Costa Rica'
$ws.Range("B84").Value = 7596
$ws.Range("D84").Value = 2239
$ws.Range("E84").Value = 5327
$ws.Range("H84").Value = 30

# Liberia <-> Republica de Chipre (rows 141/142)
$ws.Range("A141").Value = 'Liberia'
$ws.Range("B141").Value = 1024
$ws.Range("C141").Value = 14
$ws.Range("D141").Value = 439
$ws.Range("E141").Value = 534
$ws.Range("H141").Value = 51
$ws.Range("A142").Value = 'Republica de Chipre'
$ws.Range("B142").Value = 1021
$ws.Range("D142").Value = 839
$ws.Range("E142").Value = 163
$ws.Range("H142").Value = 19

# Surinam <-> Jamaica (rows 149/150)
$ws.Range("A149").Value = 'Surinam'
$ws.Range("B149").Value = 762
$ws.Range("C149").Value = 21
$ws.Range("D149").Value = 509
$ws.Range("E149").Value = 235
$ws.Range("H149").Value = 18
$ws.Range("A150").Value = 'Jamaica'
$ws.Range("B150").Value = 758
$ws.Range("D150").Value = 620
$ws.Range("E150").Value = 128
$ws.Range("H150").Value = 10

# --- Burundi jumps ahead of Martinica/Lesoto/Eritrea/Mongolia/Islas Caimanes ---
# now Burundi (row 168)
$ws.Range("A168").Value = 'Burundi'
$ws.Range("B168").Value = 269
$ws.Range("C168").Value = 78
$ws.Range("D168").Value = 207
$ws.Range("E168").Value = 61
$ws.Range("H168").Value = 1

# now Martinica (row 169)
$ws.Range("A169").Value = 'Martinica'
$ws.Range("B169").Value = 255
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 98
$ws.Range("E169").Value = 142
$ws.Range("H169").Value = 15

# now Lesoto (row 170)
$ws.Range("A170").Value = 'Lesoto'
$ws.Range("B170").Value = 245
$ws.Range("C170").Value = 12
$ws.Range("D170").Value = 32
$ws.Range("E170").Value = 211
$ws.Range("H170").Value = 2

# now Eritrea (row 171)
$ws.Range("A171").Value = 'Eritrea'
$ws.Range("B171").Value = 232
$ws.Range("D171").Value = 107
$ws.Range("E171").Value = 125

# now Mongolia (row 172)
$ws.Range("A172").Value = 'Mongolia'
$ws.Range("B172").Value = 230
$ws.Range("D172").Value = 203
$ws.Range("E172").Value = 27
$ws.Range("H172").Value = 0

# now Islas Caimanes (row 173)
$ws.Range("A173").Value = 'Islas Caimanes'
$ws.Range("B173").Value = 201
$ws.Range("D173").Value = 197
$ws.Range("E173").Value = 3

# --- Islas Malvinas <-> Groenlandia (tied case counts, order swaps) ---
$ws.Range("A209").Value = 'Islas Malvinas'
$ws.Range("A210").Value = 'Groenlandia'
